$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Cadastro de gastos com água" phase: swap energia/água content ---
# Week 3-4 (D12) used to be about "energia"; it's now about "água" and gets
# underlined emphasis (new font, still wrapped + vertically centered).
$ws.Range("D12").Value = "Melhoria HomePage; Cadastro/Atualização/Exclusão dos gastos com água; Consulta do gasto com água;                           "
$ws.Range("D12").Font.Underline = $true

# Week 5-6 (D13) used to be about "água"; it's now about "energia" (plus lixo)
$ws.Range("D13").Value = "Cadastro/Atualização/Exclusão dos gastos com energia; Consulta dos gastos com energia; Cadastro/Atualização/Exclusão do descarte de lixo; Consulta do descarte de lixo;                                   "

# --- Swap the order of the "Controle de acesso" / "Cadastro de usuarios" rows ---
$ws.Range("D21").Value = "Cadastro/Atualização/Exclusão de usuários; Consultas de usuários; Tela de Login                   "
$ws.Range("D22").Value = "Controle de acesso aos cadastros (autenticação); Gráfico do descarte de lixo                                                  "

# --- Student list (C5) now includes the two new group members ---
$ws.Range("C5").Value = "Felipe Lourenci Buniatti, Gabriela Marini Maroni e Richard Gehlen Castilhos"

# --- Project title (C6) keeps its text but becomes bold + underlined, centered ---
$ws.Range("C6").Font.Underline = $true

# --- Move the active selection to the merged project-title cell ---
$ws.Range("C6:F6").Select()
